# date-format-update
# Changes the Date_Report column (AB) text values from colon-separated
# dates (MM:DD:YYYY) to hyphen-separated dates (MM-DD-YYYY) across all
# worksheets in the workbook.

$wb = $excel.ActiveWorkbook

# Map of old text -> new text for the date strings that need updating.
$dateMap = @{
    "11:24:2025" = "11-24-2025"
    "11:20:2025" = "11-20-2025"
    "11:19:2025" = "11-19-2025"
    "11:18:2025" = "11-18-2025"
}

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count + $used.Row - 1

    for ($r = 1; $r -le $lastRow; $r++) {
        $cell = $ws.Range("AB" + $r)
        $val = $cell.Value2
        if ($null -ne $val -and $dateMap.ContainsKey([string]$val)) {
            # Force the assignment to stay plain text (otherwise Excel
            # auto-detects the hyphenated string as a date and converts
            # it to a date serial number), then restore the cell's
            # original (default/unstyled) appearance.
            $cell.NumberFormat = "@"
            $cell.Value = $dateMap[[string]$val]
            $cell.Style = "Normal"
        }
    }
}
